$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

# Row 2 (Engenharia Mecanica - Campus Arcos)
$ws.Range("E2").Value = 51
$ws.Range("F2").Value = 31
$ws.Range("H2").Value = 43

# Row 19
$ws.Range("E19").Value = 56

# Row 38
$ws.Range("E38").Value = 76

# Row 41
$ws.Range("E41").Value = 40

# Row 63
$ws.Range("E63").Value = 34

# Row 76
$ws.Range("E76").Value = 51
$ws.Range("F76").Value = 18
$ws.Range("H76").Value = 35

# Row 80
$ws.Range("E80").Value = 27

# Row 89
$ws.Range("F89").Value = 16
$ws.Range("H89").Value = 23

$wb.Save()
